# WSTA_L8_lexical_semantics.pptx update
# - Title slide: "LECTURE 4" -> "LECTURE 8"
# - "Three Kinds of semantics" slide: reword distributional-semantics bullet
# - "Basic Lexical Relations" slide: trim hyponyms/troponyms + entailment wording
# - Remove the "Text Analysis Lexicons" slide (slide 18), which pushes the
#   following slides ("Multiword Lexicons", "Moving on to the corpus",
#   "Further reading") up by one position.

$p = $ppt.ActivePresentation

# --- Slide 1: title slide lecture number ---
$titleSlide = $p.Slides.Item(1)
$subTitle = $titleSlide.Shapes.Item(2).TextFrame.TextRange
$full = $subTitle.Text
$old = "LECTURE 4 "
$idx = $full.IndexOf($old)
if ($idx -ge 0) {
    $sub = $subTitle.Characters($idx + 1, $old.Length)
    $sub.Text = "LECTURE 8 "
}

# --- Slide 3: "Three Kinds of semantics" ---
$semSlide = $p.Slides.Item(3)
$semBody = $semSlide.Shapes.Item(2).TextFrame.TextRange
$full = $semBody.Text
$old = "How the position of words in texts reflect their meaning"
$new = "How words that appear together reflect their meaning"
$idx = $full.IndexOf($old)
if ($idx -ge 0) {
    $sub = $semBody.Characters($idx + 1, $old.Length)
    $sub.Text = $new
}

# --- Slide 5: "Basic Lexical Relations" ---
$relSlide = $p.Slides.Item(5)
$relBody = $relSlide.Shapes.Item(2).TextFrame.TextRange

$full = $relBody.Text
$old = "hyponyms/troponyms (specific/manner)"
$new = "hyponyms (specific)"
$idx = $full.IndexOf($old)
if ($idx -ge 0) {
    $sub = $relBody.Characters($idx + 1, $old.Length)
    $sub.Text = $new
}

$full = $relBody.Text
$old = " (whole); entailment"
$new = " (whole)"
$idx = $full.IndexOf($old)
if ($idx -ge 0) {
    $sub = $relBody.Characters($idx + 1, $old.Length)
    $sub.Text = $new
}

# --- Remove the "Text Analysis Lexicons" slide ---
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.Shapes.Item(1).TextFrame.TextRange.Text -eq "Text Analysis Lexicons") {
        $candidate.Delete()
        break
    }
}
